$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 518
$ws.Range("I28").Value = 528.3333
$ws.Range("J28").Value = 502.5
$ws.Range("K28").Value = 528.3333
$ws.Range("L28").Value = 502.5
$ws.Range("M28").Value = -43.33330000000001
$ws.Range("N28").Value = -1472.5

$ws.Range("H41").Value = 679
$ws.Range("I41").Value = 599.6667
$ws.Range("J41").Value = 758.3333
$ws.Range("K41").Value = 599.6667
$ws.Range("L41").Value = 758.3333
$ws.Range("M41").Value = -159.6667
$ws.Range("N41").Value = -1638.3333

$ws.Range("H51").Value = 2550
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 2550
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 2550
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -3518

$ws.Range("H53").Value = 290.8
$ws.Range("I53").Value = 194
$ws.Range("J53").Value = 375.5
$ws.Range("K53").Value = 194
$ws.Range("L53").Value = 375.5
$ws.Range("M53").Value = 443
$ws.Range("N53").Value = -1649.5

$ws.Range("H98").Value = 3541.1
$ws.Range("I98").Value = 1712.3334
$ws.Range("J98").Value = 20000
$ws.Range("K98").Value = 1712.3334
$ws.Range("L98").Value = 20000
$ws.Range("M98").Value = -214.3334

$ws.Range("H107").Value = 2031.6923
$ws.Range("I107").Value = 2031.6923
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2031.6923
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -111.6922999999999

$ws.Range("H122").Value = 3541.1
$ws.Range("I122").Value = 1712.3334
$ws.Range("J122").Value = 20000
$ws.Range("K122").Value = 5137.0002
$ws.Range("L122").Value = 60000
$ws.Range("M122").Value = -2687.0002

$ws.Range("H137").Value = 1902.3684
$ws.Range("I137").Value = 1790.375
$ws.Range("J137").Value = 2499.6667
$ws.Range("K137").Value = 5371.125
$ws.Range("L137").Value = 7499.000100000001
$ws.Range("M137").Value = -2821.125

$ws.Range("H138").Value = 7918.263
$ws.Range("I138").Value = 5000
$ws.Range("J138").Value = 8080.3887
$ws.Range("K138").Value = 15000
$ws.Range("L138").Value = 24241.1661
$ws.Range("M138").Value = -9860
$ws.Range("N138").Value = -34521.1661

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 50
$ws.Range("I5").Value = 50
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 50
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 62

$ws.Range("H32").Value = 10782.369
$ws.Range("I32").Value = 8230.462
$ws.Range("J32").Value = 25000.143
$ws.Range("K32").Value = 8230.462
$ws.Range("L32").Value = 25000.143
$ws.Range("M32").Value = -7943.462

$ws.Range("H97").Value = 3002
$ws.Range("I97").Value = 536
$ws.Range("J97").Value = 7317.5
$ws.Range("K97").Value = 536
$ws.Range("L97").Value = 7317.5
$ws.Range("M97").Value = -40

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 50
$ws.Range("I4").Value = 50
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 50
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 65

$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

$ws.Range("H82").Value = 49770
$ws.Range("I82").Value = 24257
$ws.Range("J82").Value = 75283
$ws.Range("K82").Value = 24257
$ws.Range("L82").Value = 75283
$ws.Range("M82").Value = -23874
$ws.Range("N82").Value = -76049

$ws.Range("H85").Value = 49770
$ws.Range("I85").Value = 24257
$ws.Range("J85").Value = 75283
$ws.Range("K85").Value = 24257
$ws.Range("L85").Value = 75283
$ws.Range("M85").Value = -22931
$ws.Range("N85").Value = -77935

$ws.Range("H94").Value = 4970
$ws.Range("I94").Value = 5391.4287
$ws.Range("J94").Value = 3495
$ws.Range("K94").Value = 5391.4287
$ws.Range("L94").Value = 3495
$ws.Range("M94").Value = -4940.4287

$ws.Range("H132").Value = 70000
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 70000
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 70000
$ws.Range("N132").Value = -80120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2247
$ws.Range("I16").Value = 2500
$ws.Range("J16").Value = 2196.4
$ws.Range("K16").Value = 2500
$ws.Range("L16").Value = 2196.4
$ws.Range("M16").Value = -2213
$ws.Range("N16").Value = -2770.4

$ws.Range("H38").Value = 40000
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 40000
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 40000
$ws.Range("M38").ClearContents()
$ws.Range("N38").Value = -40754

$ws.Range("H46").Value = 40000
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 40000
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 40000
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -40422

$ws.Range("H86").Value = 7098.125
$ws.Range("I86").Value = 6398.3335
$ws.Range("J86").Value = 9197.5
$ws.Range("K86").Value = 6398.3335
$ws.Range("L86").Value = 9197.5
$ws.Range("M86").Value = -5275.3335

$ws.Range("H89").Value = 7098.125
$ws.Range("I89").Value = 6398.3335
$ws.Range("J89").Value = 9197.5
$ws.Range("K89").Value = 31991.6675
$ws.Range("L89").Value = 45987.5
$ws.Range("M89").Value = -26375.6675

$ws.Range("H113").Value = 2247
$ws.Range("I113").Value = 2500
$ws.Range("J113").Value = 2196.4
$ws.Range("K113").Value = 2500
$ws.Range("L113").Value = 2196.4
$ws.Range("M113").Value = -330
$ws.Range("N113").Value = -6536.4

$ws.Range("H122").Value = 2404.5908
$ws.Range("I122").Value = 2349.9412
$ws.Range("J122").Value = 2590.4
$ws.Range("K122").Value = 7049.823600000001
$ws.Range("L122").Value = 7771.200000000001
$ws.Range("M122").Value = -4599.823600000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 2140.5
$ws.Range("I22").Value = 2140.5
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 6421.5
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -6252.5

$ws.Range("H27").Value = 2140.5
$ws.Range("I27").Value = 2140.5
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 6421.5
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -6319.5

$ws.Range("H32").Value = 2000
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 2000
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 6000
$ws.Range("N32").Value = -6566

$ws.Range("H38").Value = 614.5
$ws.Range("I38").Value = 597.5
$ws.Range("J38").Value = 648.5
$ws.Range("K38").Value = 1792.5
$ws.Range("L38").Value = 1945.5
$ws.Range("M38").Value = -1445.5
$ws.Range("N38").Value = -2639.5

$ws.Range("H39").Value = 15249.5
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 15249.5
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 45748.5
$ws.Range("N39").Value = -46336.5

$ws.Range("H46").Value = 499
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 499
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 1497
$ws.Range("N46").Value = -1679

$ws.Range("H58").Value = 2700
$ws.Range("I58").Value = 2500
$ws.Range("J58").Value = 2800
$ws.Range("K58").Value = 7500
$ws.Range("L58").Value = 8400
$ws.Range("M58").Value = -7372
$ws.Range("N58").Value = -8656

$ws.Range("H98").Value = 2861.1667
$ws.Range("I98").Value = 3155.4
$ws.Range("J98").Value = 1390
$ws.Range("K98").Value = 9466.200000000001
$ws.Range("L98").Value = 4170
$ws.Range("M98").Value = -7968.200000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1677.1875
$ws.Range("I113").Value = 1630.4546
$ws.Range("J113").Value = 1780
$ws.Range("K113").Value = 1630.4546
$ws.Range("L113").Value = 1780
$ws.Range("M113").Value = 539.5454

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4446.3335
$ws.Range("I7").Value = 4446.3335
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 4446.3335
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -4334.3335
$ws.Range("N7").ClearContents()

$ws.Range("H22").Value = 6256.778
$ws.Range("I22").Value = 5157.625
$ws.Range("J22").Value = 7136.1
$ws.Range("K22").Value = 5157.625
$ws.Range("L22").Value = 7136.1
$ws.Range("M22").Value = -4862.625
$ws.Range("N22").Value = -7726.1

$ws.Range("H27").Value = 6256.778
$ws.Range("I27").Value = 5157.625
$ws.Range("J27").Value = 7136.1
$ws.Range("K27").Value = 5157.625
$ws.Range("L27").Value = 7136.1
$ws.Range("M27").Value = -5050.625
$ws.Range("N27").Value = -7350.1

$ws.Range("H126").Value = 4446.3335
$ws.Range("I126").Value = 4446.3335
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 13339.0005
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -10869.0005
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
